# EPBDS-11051 Test is added
# Reproduces the addition of a "FullJavaClassName" test table (rows 16-21)
# on the "Rules" sheet, plus the accompanying re-format of the existing
# test table (rows 3-10): borders removed, cells unified to the plain
# wrap-text / vertically-centered style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Re-format the existing "test1/test2" rows: drop the thin border
#    that used to box each value in, keep wrap text + vertical centering.
# ---------------------------------------------------------------------
$reformatCells = @("B3","D3","B4","D4","B7","D7","B8","D8","B9","D9","B10","D10")
foreach ($addr in $reformatCells) {
    $c = $ws.Range($addr)
    $c.Style = "Normal"
    $c.WrapText = $true
    $c.VerticalAlignment = -4108   # xlCenter (vertical)
}

# ---------------------------------------------------------------------
# 2. Touch every cell of the new A16:G21 block so it becomes part of the
#    sheet's used range, even the cells that stay blank.
# ---------------------------------------------------------------------
$blankCells = @(
    "A16","B16","C16","D16","E16","F16","G16",
    "A17","B17","C17","D17","E17","F17","G17",
    "A18","C18","D18","E18","F18","G18",
    "A19","D19","E19","F19","G19",
    "A20","D20","E20","F20","G20",
    "A21","D21","E21","F21","G21"
)
foreach ($addr in $blankCells) {
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# 3. Fill in the new "mapSpr" spreadsheet-method test data.
# ---------------------------------------------------------------------
$ws.Range("B18").Style = "Normal"
$ws.Range("B18").Value = "Spreadsheet SpreadsheetResult  mapSpr(String a, Map bMap, String java)"

$ws.Range("B19").Style = "Normal"
$ws.Range("B19").Value = "Steps"
$ws.Range("C19").Style = "Normal"
$ws.Range("C19").Value = "Values"

$ws.Range("B20").Style = "Normal"
$ws.Range("B20").Value = "step1"
$ws.Range("C20").Style = "Normal"
$ws.Range("C20").Value = "'= java.lang.String.class"

$ws.Range("B21").Style = "Normal"
$ws.Range("B21").Value = "step2"
$ws.Range("C21").Style = "Normal"
$ws.Range("C21").Value = "'= java.lang.Boolean.TRUE;"

# ---------------------------------------------------------------------
# 4. Match the final selection left behind in the saved file.
# ---------------------------------------------------------------------
[void]$ws.Range("C20").Select()

Write-Output "FullJavaClassName test rows added"
